$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.919.91'
$ws.Range("E2").Value = '  -2.83%  '
$ws.Range("D3").Value = '2.286.35'
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.55'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.66'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").Value = '2.284.33'
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0939'
$ws.Range("E10").Value = '  -3.59%  '
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.66'
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '2.674.76'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.55'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '53.869.00'
$ws.Range("E16").Value = '  -2.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000129'
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").Value = '2.293.62'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.90'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.03'
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '299.06'
$ws.Range("E21").Value = '  -2.77%  '
$ws.Range("E22").Value = '  +1.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("D27").Value = '2.394.37'
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.13'
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '162.88'
$ws.Range("E30").Value = '  -5.78%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = '0.0₃0682'
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.82'
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.49'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.18'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.857'
$ws.Range("E39").Value = '  +4.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.43'
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +5.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.21'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0890'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.548'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '237.98'
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0479'
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("E51").Value = '  -0.52%  '
